$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4300532569681106
$ws.Range("C2").Value = 0.6381869306858614
$ws.Range("D2").Value = 0.2030773565271132
$ws.Range("E2").Value = 0.1267143271622446
$ws.Range("F2").Value = 0.1083892220473013
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0.1040336513182604

$ws.Range("B3").Value = 0.2912603251506227
$ws.Range("C3").Value = 0.4555899203304116
$ws.Range("D3").Value = 0.500702505273356
$ws.Range("E3").Value = 0.6478556975768642
$ws.Range("F3").Value = 0.01985587814474775
$ws.Range("G3").Value = 0.01503155108693894
$ws.Range("H3").Value = 0.005031856978481451
$ws.Range("I3").Value = 0.2305284408284464

$ws.Range("B4").Value = 0.0647272092679665
$ws.Range("C4").Value = 0.08963759379215461
$ws.Range("D4").Value = 0.09848339444296765
$ws.Range("E4").Value = 0.03446581790823694
$ws.Range("F4").Value = 0.07896354873098892
$ws.Range("G4").Value = 0.08514855128525896
$ws.Range("H4").Value = 0.009206583418526912
$ws.Range("I4").Value = 0.2032957249758403

$ws.Range("B5").Value = 0.03060995832026178
$ws.Range("C5").Value = 0.4134696853263975
$ws.Range("D5").Value = 0.4173858467448348
$ws.Range("E5").Value = 0.2090704852622445
$ws.Range("F5").Value = 0.07698366186352119
$ws.Range("G5").Value = 0.07698366186352119
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = -0.0205976685375932

$ws.Range("B6").Value = 0.06755528627688583
$ws.Range("C6").Value = 0.235664965475833
$ws.Range("D6").Value = 0.1412429378531074
$ws.Range("E6").Value = 0.257045890512026
$ws.Range("F6").Value = 0.05026381560677583
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0.03056768558951965
$ws.Range("I6").Value = 0.0961587382345459

$ws.Range("B7").Value = -0.01338623961886766
$ws.Range("C7").Value = 0.09277803310310895
$ws.Range("D7").Value = -0.003494602780150706
$ws.Range("E7").Value = -0.01481582475515261
$ws.Range("F7").Value = -0.01262653256239572
$ws.Range("G7").Value = 0.0270446635730858
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = -0.0490897192811718

$ws.Range("B8").Value = 0.3192404076366557
$ws.Range("C8").Value = 0.4556137930660287
$ws.Range("D8").Value = 0.1518391450736634
$ws.Range("E8").Value = 0.2995143080226295
$ws.Range("F8").Value = -0.003810412510248724
$ws.Range("G8").Value = -0.003358557855611708
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0.03875833670253555

$ws.Range("B9").Value = 0.09079775386990885
$ws.Range("C9").Value = 0.1906086956521739
$ws.Range("D9").Value = 0.09397813753978128
$ws.Range("E9").Value = 0.1274412973524069
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0.425048904774555

$ws.Range("B10").Value = 0.1357636706473916
$ws.Range("C10").Value = -0.01508485229415466
$ws.Range("D10").Value = 0.03097252324687125
$ws.Range("E10").Value = -0.00678179402922467
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0.1645177926766375
$ws.Range("I10").Value = -0.06427090532135461

$ws.Range("B11").Value = 0.3479565574108277
$ws.Range("C11").Value = 0.3667346279398822
$ws.Range("D11").Value = 0.01673303901252364
$ws.Range("E11").Value = 0.1654908557075328
$ws.Range("F11").Value = 0.267064606741573
$ws.Range("G11").Value = 0.4361055656714126
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0.2027503702696708

$ws.Range("B12").Value = 0.3739460370994941
$ws.Range("C12").Value = 0.04454621149042463
$ws.Range("D12").Value = -0.02704987320371931
$ws.Range("E12").Value = -0.02704987320371931
$ws.Range("F12").Value = 0.2873900293255132
$ws.Range("G12").Value = 0.2873900293255132
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0.02302631578947371

$ws.Range("B13").Value = 0.4878048780487805
$ws.Range("C13").Value = 0.3567567567567567
$ws.Range("D13").Value = 0.07297297297297292
$ws.Range("E13").Value = 0.4878048780487805
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.8076923076923077
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0.1384615384615385

$ws.Range("B14").Value = 0.08835616438356159
$ws.Range("C14").Value = -0.01712328767123293
$ws.Range("D14").Value = -0.002054794520548004
$ws.Range("E14").Value = -0.03125000000000002
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = -0.1016566265060242

$ws.Range("B15").Value = 0.2361035948975647
$ws.Range("C15").Value = 0.4743097800655124
$ws.Range("D15").Value = 0.1881925522252498
$ws.Range("E15").Value = 0.3415634948708878
$ws.Range("F15").Value = 0.03383685800604227
$ws.Range("G15").Value = 0.09124805800103572
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0.2936329588014981

$ws.Range("B16").Value = 0.3347833066793463
$ws.Range("C16").Value = 0.07506651050544082
$ws.Range("D16").Value = 0.4625020015806031
$ws.Range("E16").Value = -0.006931881543836312
$ws.Range("F16").Value = 0.1015716846978238
$ws.Range("G16").Value = 0.05903039068493754
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0.06975112068370014

$ws.Range("B17").Value = 0.007513019721676731
$ws.Range("C17").Value = -0.07526881720430102
$ws.Range("D17").Value = 0.1268575266560322
$ws.Range("E17").Value = -0.07606706859833832
$ws.Range("F17").Value = -0.0683102208525939
$ws.Range("G17").Value = -0.0461798583958576
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0.06448875689435726

$ws.Range("B18").Value = 0.04347010049952089
$ws.Range("C18").Value = 0.08575866630720322
$ws.Range("D18").Value = 0.03441228604967895
$ws.Range("E18").Value = 0.03403622967768067
$ws.Range("F18").Value = 0.07379190862354047
$ws.Range("G18").Value = 0.04243496252205133
$ws.Range("H18").Value = 0.004116510690518683
$ws.Range("I18").Value = 0.1731200172646994

$ws.Range("B19").Value = 0.4253259016328846
$ws.Range("C19").Value = 0.4344161159136785
$ws.Range("D19").Value = 0.5430970743035366
$ws.Range("E19").Value = 0.3547427597884241
$ws.Range("F19").Value = -0.002065805046665021
$ws.Range("G19").Value = 0.019972916490203
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0.3584307522579702

$ws.Range("B20").Value = -0.01927710843373487
$ws.Range("C20").Value = 0.01614832535885164
$ws.Range("D20").Value = 0.1931710997924921
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0.000446162998215441
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = -0.003409933283914061

$ws.Range("B21").Value = 0.01931837524220034
$ws.Range("C21").Value = 0.01096977864089497
$ws.Range("D21").Value = -0.01489429774628533
$ws.Range("E21").Value = 0.004835692820480729
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0.01584586462328555

$ws.Range("B22").Value = 0.004263202512879448
$ws.Range("C22").Value = 0.00372855376539616
$ws.Range("D22").Value = 0.09759228247671256
$ws.Range("E22").Value = -0.003560845240907797
$ws.Range("F22").Value = -0.008272632452480092
$ws.Range("G22").Value = -0.006264478707435987
$ws.Range("H22").Value = -0.004807365766138682
$ws.Range("I22").Value = 0.00539284590237468
